# Beulah Heights University Organizations workbook restructuring:
#  - Insert a new leading "Category" column (old "Categories" column B moves to A,
#    old "Organization Name" column A moves to B).
#  - Rename several headers.
#  - Drop the (always empty) "Website" column and shift the social-link columns
#    left by one so LinkedIn/Instagram/Facebook/Twitter occupy H/I/J/K.
#  - Append two brand-new trailing columns: "Youtube Link" (L) and "Tiktok Link" (M).
#  - Resize columns to the widths used in the refreshed template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 5

# --- Step 1: remember the data that needs to move / survive the reshuffle ----
# old column A = Organization Name, old column B = Categories
$orgName = @{}
$category = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $orgName[$r]  = $ws.Cells.Item($r, 1).Value()
    $category[$r] = $ws.Cells.Item($r, 2).Value()
}

# old column H (Website) data is discarded; remember I..L (LinkedIn..Twitter) so
# they can be shifted left into H..K.
$colI = @{}
$colJ = @{}
$colK = @{}
$colL = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $colI[$r] = $ws.Cells.Item($r, 9).Value()
    $colJ[$r] = $ws.Cells.Item($r, 10).Value()
    $colK[$r] = $ws.Cells.Item($r, 11).Value()
    $colL[$r] = $ws.Cells.Item($r, 12).Value()
}

# --- Step 2: rewrite row 1 headers -------------------------------------------
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Organization Name"
$ws.Range("C1").Value = "Organization Link"
$ws.Range("D1").Value = "Logo Link"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Phone Number"
$ws.Range("H1").Value = "Linkedin Link"
$ws.Range("I1").Value = "Instagram Link"
$ws.Range("J1").Value = "Facebook Link"
$ws.Range("K1").Value = "Twitter Link"
$ws.Range("L1").Value = "Youtube Link"
$ws.Range("M1").Value = "Tiktok Link"

# New header cell M1 needs the same look as the rest of the header row
# (bold font, thin box border, centered horizontally, top vertically).
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: rewrite the data rows for the shuffled columns ------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $category[$r]
    $ws.Cells.Item($r, 2).Value = $orgName[$r]

    $ws.Cells.Item($r, 8).Value  = $colI[$r]
    $ws.Cells.Item($r, 9).Value  = $colJ[$r]
    $ws.Cells.Item($r, 10).Value = $colK[$r]
    $ws.Cells.Item($r, 11).Value = $colL[$r]

    $ws.Cells.Item($r, 12).Value = ""
    $ws.Cells.Item($r, 13).Value = ""
}

# --- Step 4: resize columns to match the refreshed template ------------------
# Excel's ColumnWidth property is offset from the OOXML stored width by the
# standard glyph-padding constant (~0.8333 characters), so subtract it off to
# land exactly on the target stored widths.
$pad = 0.8333333333333334
$ws.Columns.Item(1).ColumnWidth  = 10 - $pad
$ws.Columns.Item(2).ColumnWidth  = 19 - $pad
$ws.Columns.Item(3).ColumnWidth  = 42 - $pad
$ws.Columns.Item(4).ColumnWidth  = 50 - $pad
$ws.Columns.Item(5).ColumnWidth  = 50 - $pad
$ws.Columns.Item(6).ColumnWidth  = 27 - $pad
$ws.Columns.Item(7).ColumnWidth  = 14 - $pad
$ws.Columns.Item(8).ColumnWidth  = 15 - $pad
$ws.Columns.Item(9).ColumnWidth  = 16 - $pad
$ws.Columns.Item(10).ColumnWidth = 15 - $pad
$ws.Columns.Item(11).ColumnWidth = 14 - $pad
$ws.Columns.Item(12).ColumnWidth = 14 - $pad
$ws.Columns.Item(13).ColumnWidth = 13 - $pad
